# AutoCommit_27 октября 2023 г. 10:40:58_SibNout2023
# Fills in attendance/grade marks ("ОК"/"ок") across several rows of the
# gradebook sheet, extends a few rows with new marked cells in previously
# empty columns, and moves the frozen-pane scroll position / active cell
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$OK = "ОК"
$ok = "ок"

# ---------------------------------------------------------------------
# New cells that need the same "data cell" formatting (thick box border,
# centered, wrapped) as the rest of the grid before we can put values in
# them. We copy the format from a neighboring already-styled data cell
# (style index 2) and paste only the formatting, then set the value.
# ---------------------------------------------------------------------
$templateCell = $ws.Range("F8")
$templateCell.Copy()
$newlyFormattedCells = "G8","H8","H9","G18","H20","G25","I25","G27","H27"
foreach ($addr in $newlyFormattedCells) {
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# G8 stays blank (only gains the border/format, no value).
$ws.Range("H8").Value = $OK
$ws.Range("H9").Value = $OK
$ws.Range("G18").Value = $OK
$ws.Range("H20").Value = $OK
$ws.Range("G25").Value = $OK
$ws.Range("I25").Value = $OK
$ws.Range("G27").Value = $OK
$ws.Range("H27").Value = $OK

# ---------------------------------------------------------------------
# Existing (already styled) cells that simply get a value written in.
# ---------------------------------------------------------------------
$ws.Range("E8").Value = $ok
$ws.Range("F8").Value = $OK

$ws.Range("C9").Value = $OK

$ws.Range("F13").Value = $OK

$ws.Range("C18").Value = $OK

$ws.Range("D20").Value = $OK

$ws.Range("E23").Value = $OK

$ws.Range("C24").Value = $OK

$ws.Range("C26").Value = $OK

$ws.Range("D27").Value = $OK
$ws.Range("E27").Value = $OK
$ws.Range("F27").Value = $OK

# ---------------------------------------------------------------------
# View state: keep the freeze pane at its frozen boundary (C5, i.e. row 5
# / column C, just below/right of the frozen header rows/cols) while the
# active selection moves to C18.
# ---------------------------------------------------------------------
$ws.Range("C18").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 3
